$wb = $excel.ActiveWorkbook

# --- Large Company sheet: update HP Labs CEO/Founders info and add Microsoft Research row ---
$lc = $wb.Worksheets.Item("Large Company")

# Update existing row 22 (Hewlett Packard Labs) CEO/Founders/Affiliation cell
$lc.Cells.Item(22, 5).Value2 = "CTO Tolga Kurtoglu; Fabian Böhm,, Bassem Tossoun, Thomas Van Vaerenbergh"

# Add a brand new company row (row 23): Microsoft Research
$lc.Cells.Item(23, 2).Value2 = "Microsoft Research"
$lc.Cells.Item(23, 3).Value2 = 100
$lc.Cells.Item(23, 4).Value2 = "Analog Optical Computing. MicroLed Array, Free Space MVM, Electrical Nonlinearity"
$lc.Cells.Item(23, 5).Value2 = "H. Ballani, G. Brennan, B. Canakci, J. Chu, J. H. Clegg, D. Cletheroe, C. Gkantsidis, J. Gladrow, K. P.`nKalinin, D. J. Kelly, H. Kremer, G. O'Shea, F. Parmigiani, L. Pickup, B. Rahmani, A. Rowstron"
$lc.Cells.Item(23, 6).Value2 = "Cambridge, England"

# Make "Large Company" the active/selected sheet, scrolled/selected near the newly added row
$lc.Activate()
$lc.Range("F19").Select()
$excel.ActiveWindow.ScrollRow = 6
$excel.ActiveWindow.ScrollColumn = 1
